$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rng = $ws.Range("A2:D90")
$key = $ws.Range("A1:A90")

$rng.Sort($key, 1)
